# Adds a new "14-10-2020" column (AC) to the COVID19_TIMESERIESDATA sheet:
#   - AC1 header gets the same look as the other recent date headers
#     (bold / centred / thin-bordered, matching AB1's formatting)
#   - AC2:AC36 get the per-state case counts for 2020-10-14
# This extends the sheet's used range from A1:AB36 to A1:AC36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID19_TIMESERIESDATA")

# --- Header cell AC1 -------------------------------------------------
# Copy AB1's formatting (border/bold font/centred alignment) onto AC1,
# then set its own text value.
$ws.Range("AB1").Copy()
$ws.Range("AC1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("AC1").Value = "14-10-2020"

# --- Data values for column AC (2020-10-14), rows 2-36 ---------------
$values = @(
    3782,
    714427,
    9573,
    167059,
    187059,
    12007,
    119352,
    3081,
    286880,
    34252,
    134990,
    132382,
    15001,
    73502,
    85314,
    602505,
    207357,
    4205,
    132429,
    1297252,
    10829,
    5406,
    2093,
    5831,
    232988,
    26865,
    113105,
    139616,
    3014,
    612320,
    191269,
    25041,
    48283,
    397570,
    265288
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 29).Value = $values[$i]
}
